$d = $word.ActiveDocument

# 1) Merge "My 20 " / "years experience" / " is primarily..." into one run (remove proofErr spell-check split)
$d.Content.Find.Execute("My 20 years experience is primarily in the field of data warehousing and ETL technology.", $false, $false, $false, $false, $false, $true, 1, $false, "My 20 years experience is primarily in the field of data warehousing and ETL technology.", 2) | Out-Null

# 2) Merge "This was my specialization in the " / "Masters" / " Degree..." into one run
$d.Content.Find.Execute("This was my specialization in the Masters Degree Computer Science that I completed from Villanova University, Philadelphia USA.", $false, $false, $false, $false, $false, $true, 1, $false, "This was my specialization in the Masters Degree Computer Science that I completed from Villanova University, Philadelphia USA.", 2) | Out-Null

# 3) Merge "My experience is in the ETL tool IBM " / "Datastage" into one run (do not touch following runs)
$d.Content.Find.Execute("My experience is in the ETL tool IBM Datastage", $false, $false, $false, $false, $false, $true, 1, $false, "My experience is in the ETL tool IBM Datastage", 2) | Out-Null

# 4) Append new sentence after "...oracle, db2, Netezza among others."
$d.Content.Find.Execute("oracle, db2, Netezza among others.", $false, $false, $false, $false, $false, $true, 1, $false, "oracle, db2, Netezza among others. In addition to this I have an ancillary ETL skill in Microsoft SSIS.", 2) | Out-Null

# 5) Merge "20+ years of development experience in the IBM stack " / "Datastage" / ". This includes..." into one run
$d.Content.Find.Execute("20+ years of development experience in the IBM stack Datastage. This includes development of generic ETL jobs for extract, transforming and loading data. ", $false, $false, $false, $false, $false, $true, 1, $false, "20+ years of development experience in the IBM stack Datastage. This includes development of generic ETL jobs for extract, transforming and loading data. ", 2) | Out-Null
